# Applies the attendance-report sync edit:
#  - Attendance sheet: insert 3 new check-in rows (two new B1-3 students'
#    first session plus a repeat-scan row) at row 435, shifting the
#    existing B1-6 "session" rows down by 3 (435-481 -> 438-484).
#  - Summary sheet: refresh the aggregated attendance stats for the three
#    students (211410, 211439, 211446) who picked up the extra session.
#  - Keep dimension / autofilter / the hidden _FilterDatabase defined name
#    in sync with the new row count (481 -> 484).

$wb = $excel.ActiveWorkbook
$attendance = $wb.Worksheets.Item("Attendance")
$summary = $wb.Worksheets.Item("Summary")

# ---------------------------------------------------------------------
# 1) Attendance sheet: insert 3 blank rows before the old row 435, then
#    populate them with the new check-ins.
# ---------------------------------------------------------------------
$attendance.Range("A435:A437").EntireRow.Insert()

# Force these brand-new cells to be stored as TEXT (matching every other
# cell in this sheet, which are all inlineStr) instead of letting the
# host auto-detect numbers/dates for values like "211439" or "08/12/2025".
$attendance.Range("A435:K437").NumberFormat = "@"

$newRows = @(
    @("211439", "عبد الرحمن سامح عبد العزيز منصور", "Year 5", "B1-3", "211439@med.asu.edu.eg", "GENERAL SURGERY", "2", "GENERAL SURGERY", "08/12/2025", "12:30:00", "B1-3"),
    @("211446", "ساره عبد الله محمد كمال عبد العزيز", "Year 5", "B1-3", "211446@med.asu.edu.eg", "GENERAL SURGERY", "2", "GENERAL SURGERY", "08/12/2025", "12:30:00", "B1-3"),
    @("211410", "جون مجدى ميخائيل سدراك", "Year 5", "B1-3", "211410@med.asu.edu.eg", "GENERAL SURGERY", "2", "GENERAL SURGERY", "08/12/2025", "12:30:00", "B1-3")
)

$cols = @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J", "K")
$startRow = 435
for ($r = 0; $r -lt $newRows.Length; $r++) {
    $rowData = $newRows[$r]
    $rowNum = $startRow + $r
    for ($c = 0; $c -lt $cols.Length; $c++) {
        $addr = $cols[$c] + $rowNum
        $attendance.Range($addr).Value = $rowData[$c]
    }
}

# Keep the sheet's autofilter covering the new used range.
$attendance.AutoFilterMode = $false
$attendance.Range("A1:K484").AutoFilter()

# ---------------------------------------------------------------------
# 2) Summary sheet: the three students who now attended GENERAL SURGERY
#    session 2 (on 08/12/2025) move from 1/15 to 2/15 sessions attended.
# ---------------------------------------------------------------------
$summaryRows = @(140, 159, 164)
foreach ($row in $summaryRows) {
    $gCell = $summary.Range("G$row")
    $gCell.NumberFormat = "@"
    $gCell.Value = "13.3%"

    $summary.Range("I$row").Value = 10
    $summary.Range("N$row").Value = 2
    $summary.Range("O$row").Value = 0
    $summary.Range("Q$row").Value = 2
    $summary.Range("S$row").Value = 1
}

# ---------------------------------------------------------------------
# 3) Keep the hidden _xlnm._FilterDatabase defined name for Attendance in
#    sync with the new last row (481 -> 484).
# ---------------------------------------------------------------------
$names = $wb.Names
for ($i = 1; $i -le $names.Count(); $i++) {
    $n = $names.Item($i)
    $nm = $n.Name()
    if ($nm -eq "Attendance!_FilterDatabase") {
        $n.RefersTo = "='Attendance'!`$A`$1:`$K`$484"
    }
}

"Edit complete"
